$wb = $excel.ActiveWorkbook

# --- Week2 sheet: drop the topLeftCell scroll-freeze, keep selection at F6 ---
$ws2 = $wb.Worksheets.Item("Week2")
$ws2.Range("F6").Select()

# --- Week3 sheet: the main content sheet for this edit ---
$ws3 = $wb.Worksheets.Item("Week3")

# New shared-string text, entered in the exact order needed so the
# sharedStrings.xml table is appended in the same sequence as the target file.
$ws3.Range("D3").Value = 'Did document write up with further suggestions + actions that could be needed to take'
$ws3.Range("D4").Value = 'Went to class - Did small write up for assignment1
Did weekly for for Mike '
$ws3.Range("D5").Value = 'Have a path to go down, looking forward to another meeting'
$ws3.Range("D7").Value = 'Filling in time between meetings with people. Managed to fill time with Ara work and did not put time towards project'
$ws3.Range("D6").Value = 'At least have a next meeting time'
$ws3.Range("E3").Value = 'Had meeting with John, Craig and Lorna. Integrating this whole project into the overall scope of tourplans strategy.
Started looking at how the over all structure will fit into the rest of a future project'
$ws3.Range("B13").Value = 'Weekly Mike meeting:'
$ws3.Range("C14").Value = 'positives out of 200 and 300 level courses'
$ws3.Range("B14").Value = '200/300 Review'
$ws3.Range("C15").Value = 'Took something from the couse and applied it'
$ws3.Range("C16").Value = 'What I learnt outside of the course that would have been helpful to teach in Ara'
$ws3.Range("B17").Value = 'Methadology'
$ws3.Range("C17").Value = 'Sumarise the theory(10 point that can be made into 100 word paragraph(1000 words total))(Learning point could be model, view, controller etc)'
$ws3.Range("C18").Value = '1000 Words on practiced theory what the industry does( can also ask what jonty did)'
$ws3.Range("C19").Value = 'How it works in the real world. Do the techniques work? Are they applied? '
$ws3.Range("E4").Value = 'Meeting with mike + had project class tutorial on QA '
$ws3.Range("E5").Value = 'Completed meeting'
$ws3.Range("E6").Value = 'Start write up of how solution will fit into the overall schema. '
$ws3.Range("E7").Value = 'Just getting my head around how the whole system will work.'
$ws3.Range("F3").Value = 'Created a project overview + Visio drawing for John. Wil get feedback on style or structure of it.'
$ws3.Range("F5").Value = 'Finished write up for John
Started thinking about new questions'
$ws3.Range("F6").Value = 'More Qs'
$ws3.Range("F7").Value = 'Understanding how different parts interact and the point of seperating all these systems'

# Column widths for the newly used D:F columns (target widths 13.28515625 /
# 23.28515625 / 19.7109375 chars; the inputs below are chosen so the
# engine's internal pixel-rounded ColumnWidth lands as close as possible
# to those target values)
$ws3.Columns.Item(4).ColumnWidth = 12.5
$ws3.Columns.Item(5).ColumnWidth = 22.5
$ws3.Columns.Item(6).ColumnWidth = 18.833333333333332

# Row heights (grown to fit the newly entered wrapped text)
$ws3.Rows.Item(3).RowHeight = 157.5
$ws3.Rows.Item(4).RowHeight = 102
$ws3.Rows.Item(5).RowHeight = 90.75
$ws3.Rows.Item(7).RowHeight = 150.75

# Re-activate Week3 (the workbook's active tab) and restore its selection
$ws3.Activate()
$ws3.Range("I4").Select()
